$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$urls = @(
    'https://www.mercadolivre.com.br/kit-canetas-gel-48-cores-glitter-metalico-pastel-e-neon-exterior-preto/p/MLB24557984#searchVariation=MLB24557984&position=2&search_layout=grid&type=product&tracking_id=99ea6a2f-141a-4527-8f23-eb37cfbe5738',
    'https://www.mercadolivre.com.br/15-canetas-fofas-variadas-atacado-papelaria-escolar-tinta-preto/p/MLB35822093#searchVariation=MLB35822093&position=3&search_layout=grid&type=product&tracking_id=99ea6a2f-141a-4527-8f23-eb37cfbe5738',
    'https://www.mercadolivre.com.br/caneta-esferografica-compactor-caixa-com-100-unidades-cor-da-tinta-azul-cor-do-exterior-cristal/p/MLB22766681#searchVariation=MLB22766681&position=4&search_layout=grid&type=product&tracking_id=99ea6a2f-141a-4527-8f23-eb37cfbe5738',
    'https://www.mercadolivre.com.br/caneta-bic-cristal-preciso-08mm-azul-caixa-c-50-cor-do-exterior-transparente/p/MLB22706525#searchVariation=MLB22706525&position=6&search_layout=grid&type=product&tracking_id=99ea6a2f-141a-4527-8f23-eb37cfbe5738',
    'https://www.mercadolivre.com.br/marcador-retroprojetor-cd-20mm-jocar-preta-cx12/p/MLB24732135#searchVariation=MLB24732135&position=8&search_layout=grid&type=product&tracking_id=99ea6a2f-141a-4527-8f23-eb37cfbe5738',
    'https://www.mercadolivre.com.br/kit-canetas-gel-48-cores-glitter-metalico-pastel-e-neon-exterior-preto/p/MLB24557984#searchVariation=MLB24557984&position=3&search_layout=grid&type=product&tracking_id=f388956c-6f7b-4731-9720-0a2a59b89c52',
    'https://www.mercadolivre.com.br/15-canetas-fofas-variadas-atacado-papelaria-escolar-tinta-preto/p/MLB35822093#searchVariation=MLB35822093&position=6&search_layout=grid&type=product&tracking_id=f388956c-6f7b-4731-9720-0a2a59b89c52',
    'https://www.mercadolivre.com.br/caneta-esferografica-compactor-caixa-com-100-unidades-cor-da-tinta-azul-cor-do-exterior-cristal/p/MLB22766681#searchVariation=MLB22766681&position=13&search_layout=grid&type=product&tracking_id=f388956c-6f7b-4731-9720-0a2a59b89c52',
    'https://www.mercadolivre.com.br/caneta-bic-cristal-preciso-08mm-azul-caixa-c-50-cor-do-exterior-transparente/p/MLB22706525#searchVariation=MLB22706525&position=18&search_layout=grid&type=product&tracking_id=f388956c-6f7b-4731-9720-0a2a59b89c52',
    'https://www.mercadolivre.com.br/marcador-retroprojetor-cd-20mm-jocar-preta-cx12/p/MLB24732135#searchVariation=MLB24732135&position=17&search_layout=grid&type=product&tracking_id=f388956c-6f7b-4731-9720-0a2a59b89c52',
    'https://www.mercadolivre.com.br/kit-canetas-gel-48-cores-glitter-metalico-pastel-e-neon-exterior-preto/p/MLB24557984#searchVariation=MLB24557984&position=2&search_layout=grid&type=product&tracking_id=aa23c399-e12c-417b-a8a6-453afd5fbb86',
    'https://www.mercadolivre.com.br/15-canetas-fofas-variadas-atacado-papelaria-escolar-tinta-preto/p/MLB35822093#searchVariation=MLB35822093&position=3&search_layout=grid&type=product&tracking_id=aa23c399-e12c-417b-a8a6-453afd5fbb86',
    'https://www.mercadolivre.com.br/caneta-esferografica-compactor-caixa-com-100-unidades-cor-da-tinta-azul-cor-do-exterior-cristal/p/MLB22766681#searchVariation=MLB22766681&position=4&search_layout=grid&type=product&tracking_id=aa23c399-e12c-417b-a8a6-453afd5fbb86',
    'https://www.mercadolivre.com.br/caneta-bic-cristal-preciso-08mm-azul-caixa-c-50-cor-do-exterior-transparente/p/MLB22706525#searchVariation=MLB22706525&position=5&search_layout=grid&type=product&tracking_id=aa23c399-e12c-417b-a8a6-453afd5fbb86',
    'https://www.mercadolivre.com.br/marcador-retroprojetor-cd-20mm-jocar-preta-cx12/p/MLB24732135#searchVariation=MLB24732135&position=6&search_layout=grid&type=product&tracking_id=aa23c399-e12c-417b-a8a6-453afd5fbb86',
    'https://www.mercadolivre.com.br/kit-canetas-gel-48-cores-glitter-metalico-pastel-e-neon-exterior-preto/p/MLB24557984#searchVariation=MLB24557984&position=3&search_layout=grid&type=product&tracking_id=98004e2a-881f-4d7a-8410-df5355e4ccf9',
    'https://www.mercadolivre.com.br/15-canetas-fofas-variadas-atacado-papelaria-escolar-tinta-preto/p/MLB35822093#searchVariation=MLB35822093&position=5&search_layout=grid&type=product&tracking_id=98004e2a-881f-4d7a-8410-df5355e4ccf9',
    'https://www.mercadolivre.com.br/caneta-esferografica-compactor-caixa-com-100-unidades-cor-da-tinta-azul-cor-do-exterior-cristal/p/MLB22766681#searchVariation=MLB22766681&position=7&search_layout=grid&type=product&tracking_id=98004e2a-881f-4d7a-8410-df5355e4ccf9',
    'https://www.mercadolivre.com.br/caneta-bic-cristal-preciso-08mm-azul-caixa-c-50-cor-do-exterior-transparente/p/MLB22706525#searchVariation=MLB22706525&position=9&search_layout=grid&type=product&tracking_id=98004e2a-881f-4d7a-8410-df5355e4ccf9',
    'https://www.mercadolivre.com.br/marcador-retroprojetor-cd-20mm-jocar-preta-cx12/p/MLB24732135#searchVariation=MLB24732135&position=17&search_layout=grid&type=product&tracking_id=98004e2a-881f-4d7a-8410-df5355e4ccf9'
)

$startRow = 443
for ($i = 0; $i -lt $urls.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $urls[$i]
}

